$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$cD = $ws.Range("D2")
$cD.NumberFormat = "@"
$cD.Value = "30.706.85"
$cD.ClearFormats()
$cE = $ws.Range("E2")
$cE.NumberFormat = "@"
$cE.Value = "  +0.63%  "
$cE.ClearFormats()

# Row 3 - Ethereum
$cD = $ws.Range("D3")
$cD.NumberFormat = "@"
$cD.Value = "1.965.74"
$cD.ClearFormats()
$cE = $ws.Range("E3")
$cE.NumberFormat = "@"
$cE.Value = "  +2.60%  "
$cE.ClearFormats()

# Row 4 - TetherUSD (price unchanged)
$cE = $ws.Range("E4")
$cE.NumberFormat = "@"
$cE.Value = "  +0.00%  "
$cE.ClearFormats()

# Row 5 - BNB
$cD = $ws.Range("D5")
$cD.NumberFormat = "@"
$cD.Value = "249.62"
$cD.ClearFormats()
$cE = $ws.Range("E5")
$cE.NumberFormat = "@"
$cE.Value = "  +1.65%  "
$cE.ClearFormats()

# Row 6 - USDC (price unchanged)
$cE = $ws.Range("E6")
$cE.NumberFormat = "@"
$cE.Value = "  +0.10%  "
$cE.ClearFormats()

# Row 7 - XRP
$cD = $ws.Range("D7")
$cD.NumberFormat = "@"
$cD.Value = "0.4834"
$cD.ClearFormats()
$cE = $ws.Range("E7")
$cE.NumberFormat = "@"
$cE.Value = "  +0.24%  "
$cE.ClearFormats()

# Row 8 - OKB
$cD = $ws.Range("D8")
$cD.NumberFormat = "@"
$cD.Value = "44.74"
$cD.ClearFormats()
$cE = $ws.Range("E8")
$cE.NumberFormat = "@"
$cE.Value = "  +0.87%  "
$cE.ClearFormats()

# Row 9 - Cardano
$cD = $ws.Range("D9")
$cD.NumberFormat = "@"
$cD.Value = "0.2954"
$cD.ClearFormats()
$cE = $ws.Range("E9")
$cE.NumberFormat = "@"
$cE.Value = "  +2.19%  "
$cE.ClearFormats()

# Row 10 - Dogecoin
$cD = $ws.Range("D10")
$cD.NumberFormat = "@"
$cD.Value = "0.06814"
$cD.ClearFormats()
$cE = $ws.Range("E10")
$cE.NumberFormat = "@"
$cE.Value = "  +1.35%  "
$cE.ClearFormats()

# Row 11 - Litecoin
$cD = $ws.Range("D11")
$cD.NumberFormat = "@"
$cD.Value = "111.27"
$cD.ClearFormats()
$cE = $ws.Range("E11")
$cE.NumberFormat = "@"
$cE.Value = "  +0.81%  "
$cE.ClearFormats()

# Row 12 - Solana
$cD = $ws.Range("D12")
$cD.NumberFormat = "@"
$cD.Value = "19.46"
$cD.ClearFormats()
$cE = $ws.Range("E12")
$cE.NumberFormat = "@"
$cE.Value = "  +1.03%  "
$cE.ClearFormats()

# Row 13 - WrappedEther
$cD = $ws.Range("D13")
$cD.NumberFormat = "@"
$cD.Value = "1.959.00"
$cD.ClearFormats()
$cE = $ws.Range("E13")
$cE.NumberFormat = "@"
$cE.Value = "  +2.24%  "
$cE.ClearFormats()

# Row 14 - TRON
$cD = $ws.Range("D14")
$cD.NumberFormat = "@"
$cD.Value = "0.07746"
$cD.ClearFormats()
$cE = $ws.Range("E14")
$cE.NumberFormat = "@"
$cE.Value = "  +2.41%  "
$cE.ClearFormats()

# Row 15 - Polkadot
$cD = $ws.Range("D15")
$cD.NumberFormat = "@"
$cD.Value = "5.501"
$cD.ClearFormats()
$cE = $ws.Range("E15")
$cE.NumberFormat = "@"
$cE.Value = "  +4.69%  "
$cE.ClearFormats()

# Row 16 - Polygon
$cD = $ws.Range("D16")
$cD.NumberFormat = "@"
$cD.Value = "0.6944"
$cD.ClearFormats()
$cE = $ws.Range("E16")
$cE.NumberFormat = "@"
$cE.Value = "  +3.63%  "
$cE.ClearFormats()

# Row 17 - BitcoinCash
$cD = $ws.Range("D17")
$cD.NumberFormat = "@"
$cD.Value = "294.61"
$cD.ClearFormats()
$cE = $ws.Range("E17")
$cE.NumberFormat = "@"
$cE.Value = "  +2.13%  "
$cE.ClearFormats()

# Row 18 - WrappedBTC
$cD = $ws.Range("D18")
$cD.NumberFormat = "@"
$cD.Value = "30.700.88"
$cD.ClearFormats()
$cE = $ws.Range("E18")
$cE.NumberFormat = "@"
$cE.Value = "  +0.58%  "
$cE.ClearFormats()

# Row 19 - Avalanche
$cD = $ws.Range("D19")
$cD.NumberFormat = "@"
$cD.Value = "13.37"
$cD.ClearFormats()
$cE = $ws.Range("E19")
$cE.NumberFormat = "@"
$cE.Value = "  +3.61%  "
$cE.ClearFormats()

# Row 20 - Uniswap
$cD = $ws.Range("D20")
$cD.NumberFormat = "@"
$cD.Value = "5.673"
$cD.ClearFormats()
$cE = $ws.Range("E20")
$cE.NumberFormat = "@"
$cE.Value = "  +3.32%  "
$cE.ClearFormats()

# Row 21 - ShibaInu
$cD = $ws.Range("D21")
$cD.NumberFormat = "@"
$cD.Value = "0.000007712"
$cD.ClearFormats()
$cE = $ws.Range("E21")
$cE.NumberFormat = "@"
$cE.Value = "  +1.44%  "
$cE.ClearFormats()

# Row 22 - WrappedliquidstakedEther2.0
$cD = $ws.Range("D22")
$cD.NumberFormat = "@"
$cD.Value = "2.213.18"
$cD.ClearFormats()
$cE = $ws.Range("E22")
$cE.NumberFormat = "@"
$cE.Value = "  +2.17%  "
$cE.ClearFormats()

# Row 23 - Dai
$cD = $ws.Range("D23")
$cD.NumberFormat = "@"
$cD.Value = "1.002"
$cD.ClearFormats()
$cE = $ws.Range("E23")
$cE.NumberFormat = "@"
$cE.Value = "  +0.19%  "
$cE.ClearFormats()

# Row 24 - BinanceUSD
$cD = $ws.Range("D24")
$cD.NumberFormat = "@"
$cD.Value = "1.003"
$cD.ClearFormats()
$cE = $ws.Range("E24")
$cE.NumberFormat = "@"
$cE.Value = "  +0.23%  "
$cE.ClearFormats()

# Row 25 - Chainlink
$cD = $ws.Range("D25")
$cD.NumberFormat = "@"
$cD.Value = "6.662"
$cD.ClearFormats()
$cE = $ws.Range("E25")
$cE.NumberFormat = "@"
$cE.Value = "  +4.18%  "
$cE.ClearFormats()

# Row 26 - Cosmos
$cD = $ws.Range("D26")
$cD.NumberFormat = "@"
$cD.Value = "9.878"
$cD.ClearFormats()
$cE = $ws.Range("E26")
$cE.NumberFormat = "@"
$cE.Value = "  +4.38%  "
$cE.ClearFormats()

# Row 27 - Monero
$cD = $ws.Range("D27")
$cD.NumberFormat = "@"
$cD.Value = "169.94"
$cD.ClearFormats()
$cE = $ws.Range("E27")
$cE.NumberFormat = "@"
$cE.Value = "  +3.28%  "
$cE.ClearFormats()

# Row 28 - EthereumClassic
$cD = $ws.Range("D28")
$cD.NumberFormat = "@"
$cD.Value = "20.22"
$cD.ClearFormats()
$cE = $ws.Range("E28")
$cE.NumberFormat = "@"
$cE.Value = "  -0.56%  "
$cE.ClearFormats()

# Row 29 - LidoDAOToken
$cD = $ws.Range("D29")
$cD.NumberFormat = "@"
$cD.Value = "2.211"
$cD.ClearFormats()
$cE = $ws.Range("E29")
$cE.NumberFormat = "@"
$cE.Value = "  +3.38%  "
$cE.ClearFormats()

# Row 30 - Stellar
$cD = $ws.Range("D30")
$cD.NumberFormat = "@"
$cD.Value = "0.1078"
$cD.ClearFormats()
$cE = $ws.Range("E30")
$cE.NumberFormat = "@"
$cE.Value = "  +1.36%  "
$cE.ClearFormats()

# Row 31 - Toncoin
$cD = $ws.Range("D31")
$cD.NumberFormat = "@"
$cD.Value = "1.442"
$cD.ClearFormats()
$cE = $ws.Range("E31")
$cE.NumberFormat = "@"
$cE.Value = "  +2.37%  "
$cE.ClearFormats()

# Row 32 - Filecoin
$cD = $ws.Range("D32")
$cD.NumberFormat = "@"
$cD.Value = "4.671"
$cD.ClearFormats()
$cE = $ws.Range("E32")
$cE.NumberFormat = "@"
$cE.Value = "  +15.97%  "
$cE.ClearFormats()

# Row 33 - InternetComputer(DFINITY)
$cD = $ws.Range("D33")
$cD.NumberFormat = "@"
$cD.Value = "4.457"
$cD.ClearFormats()
$cE = $ws.Range("E33")
$cE.NumberFormat = "@"
$cE.Value = "  +7.15%  "
$cE.ClearFormats()

# Row 34 - Hedera
$cD = $ws.Range("D34")
$cD.NumberFormat = "@"
$cD.Value = "0.05110"
$cD.ClearFormats()
$cE = $ws.Range("E34")
$cE.NumberFormat = "@"
$cE.Value = "  +2.34%  "
$cE.ClearFormats()

# Row 35 - ImmutableX
$cD = $ws.Range("D35")
$cD.NumberFormat = "@"
$cD.Value = "0.7818"
$cD.ClearFormats()
$cE = $ws.Range("E35")
$cE.NumberFormat = "@"
$cE.Value = "  +7.21%  "
$cE.ClearFormats()

# Row 36 - ARBITRUM
$cD = $ws.Range("D36")
$cD.NumberFormat = "@"
$cD.Value = "1.183"
$cD.ClearFormats()
$cE = $ws.Range("E36")
$cE.NumberFormat = "@"
$cE.Value = "  +4.36%  "
$cE.ClearFormats()

# Row 37 - VeChain
$cD = $ws.Range("D37")
$cD.NumberFormat = "@"
$cD.Value = "0.02066"
$cD.ClearFormats()
$cE = $ws.Range("E37")
$cE.NumberFormat = "@"
$cE.Value = "  +0.64%  "
$cE.ClearFormats()

# Row 38 - HuobiToken
$cD = $ws.Range("D38")
$cD.NumberFormat = "@"
$cD.Value = "2.738"
$cD.ClearFormats()
$cE = $ws.Range("E38")
$cE.NumberFormat = "@"
$cE.Value = "  +0.03%  "
$cE.ClearFormats()

# Row 39 - MXToken
$cD = $ws.Range("D39")
$cD.NumberFormat = "@"
$cD.Value = "2.716"
$cD.ClearFormats()
$cE = $ws.Range("E39")
$cE.NumberFormat = "@"
$cE.Value = "  +1.71%  "
$cE.ClearFormats()

# Row 40 - RenderToken
$cD = $ws.Range("D40")
$cD.NumberFormat = "@"
$cD.Value = "2.079"
$cD.ClearFormats()
$cE = $ws.Range("E40")
$cE.NumberFormat = "@"
$cE.Value = "  +3.27%  "
$cE.ClearFormats()

# Row 41 - Quant
$cD = $ws.Range("D41")
$cD.NumberFormat = "@"
$cD.Value = "111.74"
$cD.ClearFormats()
$cE = $ws.Range("E41")
$cE.NumberFormat = "@"
$cE.Value = "  +0.89%  "
$cE.ClearFormats()

# Row 42 - FraxShare
$cD = $ws.Range("D42")
$cD.NumberFormat = "@"
$cD.Value = "6.105"
$cD.ClearFormats()
$cE = $ws.Range("E42")
$cE.NumberFormat = "@"
$cE.Value = "  +3.67%  "
$cE.ClearFormats()

# Row 43 - TheSandbox
$cD = $ws.Range("D43")
$cD.NumberFormat = "@"
$cD.Value = "0.4484"
$cD.ClearFormats()
$cE = $ws.Range("E43")
$cE.NumberFormat = "@"
$cE.Value = "  +1.25%  "
$cE.ClearFormats()

# Row 44 - TrustWalletToken
$cD = $ws.Range("D44")
$cD.NumberFormat = "@"
$cD.Value = "0.8770"
$cD.ClearFormats()
$cE = $ws.Range("E44")
$cE.NumberFormat = "@"
$cE.Value = "  +1.48%  "
$cE.ClearFormats()

# Row 45 - Aave
$cD = $ws.Range("D45")
$cD.NumberFormat = "@"
$cD.Value = "70.28"
$cD.ClearFormats()
$cE = $ws.Range("E45")
$cE.NumberFormat = "@"
$cE.Value = "  +3.15%  "
$cE.ClearFormats()

# Row 46 - PaxDollar
$cD = $ws.Range("D46")
$cD.NumberFormat = "@"
$cD.Value = "1.003"
$cD.ClearFormats()
$cE = $ws.Range("E46")
$cE.NumberFormat = "@"
$cE.Value = "  +0.22%  "
$cE.ClearFormats()

# Row 47 - Aptos
$cD = $ws.Range("D47")
$cD.NumberFormat = "@"
$cD.Value = "7.464"
$cD.ClearFormats()
$cE = $ws.Range("E47")
$cE.NumberFormat = "@"
$cE.Value = "  +1.58%  "
$cE.ClearFormats()

# Row 48 & 49 swap: Algorand moves to row 48, EnergySwap moves to row 49,
# each with updated price/volume figures.
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$cD = $ws.Range("D48")
$cD.NumberFormat = "@"
$cD.Value = "0.1285"
$cD.ClearFormats()
$cE = $ws.Range("E48")
$cE.NumberFormat = "@"
$cE.Value = "  +3.51%  "
$cE.ClearFormats()

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$cD = $ws.Range("D49")
$cD.NumberFormat = "@"
$cD.Value = "9.379"
$cD.ClearFormats()
$cE = $ws.Range("E49")
$cE.NumberFormat = "@"
$cE.Value = "  +0.44%  "
$cE.ClearFormats()

# Row 50 - Elrond
$cD = $ws.Range("D50")
$cD.NumberFormat = "@"
$cD.Value = "35.98"
$cD.ClearFormats()
$cE = $ws.Range("E50")
$cE.NumberFormat = "@"
$cE.Value = "  +3.26%  "
$cE.ClearFormats()

# Row 51 - BitcoinSV
$cD = $ws.Range("D51")
$cD.NumberFormat = "@"
$cD.Value = "47.88"
$cD.ClearFormats()
$cE = $ws.Range("E51")
$cE.NumberFormat = "@"
$cE.Value = "  -2.38%  "
$cE.ClearFormats()
